# Update profit-calculation figures across the Sheets workbook (scheduled runner refresh).
# Values below replace stale market-board snapshots (H/I/J/K/L) and the derived
# profit deltas (M/N) for the affected rows on each class sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 599.8333
$ws.Range("I2").Value = 620
$ws.Range("K2").Value = 620
$ws.Range("M2").Value = -507

$ws.Range("H40").Value = 9500
$ws.Range("I40").Value = 9500
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 9500
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -9325
$ws.Range("N40").ClearContents()

$ws.Range("H76").Value = 6712.4116
$ws.Range("I76").Value = 6010.091
$ws.Range("J76").Value = 8000
$ws.Range("K76").Value = 6010.091
$ws.Range("L76").Value = 8000
$ws.Range("M76").Value = -5695.091
$ws.Range("N76").Value = -8630

$ws.Range("H79").Value = 6712.4116
$ws.Range("I79").Value = 6010.091
$ws.Range("J79").Value = 8000
$ws.Range("K79").Value = 6010.091
$ws.Range("L79").Value = 8000
$ws.Range("M79").Value = -4918.091
$ws.Range("N79").Value = -10184

$ws.Range("H94").Value = 4499.5
$ws.Range("J94").Value = 3000
$ws.Range("L94").Value = 3000
$ws.Range("N94").Value = -3902

$ws.Range("H132").Value = 3652.375
$ws.Range("I132").Value = 3593.7827
$ws.Range("K132").Value = 10781.3481
$ws.Range("M132").Value = -8251.348100000001

$ws.Range("H135").Value = 7353470
$ws.Range("I135").Value = 466.26666
$ws.Range("K135").Value = 4196.39994
$ws.Range("M135").Value = -1661.39994

$ws.Range("H137").Value = 6366.9644
$ws.Range("I137").Value = 2095.5789
$ws.Range("K137").Value = 6286.736699999999
$ws.Range("M137").Value = -3736.736699999999

$ws.Range("H138").Value = 6726.1626
$ws.Range("I138").Value = 1337
$ws.Range("J138").Value = 10249.846
$ws.Range("K138").Value = 4011
$ws.Range("L138").Value = 30749.538
$ws.Range("M138").Value = 1129
$ws.Range("N138").Value = -41029.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 4250
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 4250
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 4250
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -4596

$ws.Range("H32").Value = 10418888
$ws.Range("I32").Value = 10990383
$ws.Range("K32").Value = 10990383
$ws.Range("M32").Value = -10990096

$ws.Range("H41").Value = 11249.75
$ws.Range("I41").Value = 11249.75
$ws.Range("K41").Value = 11249.75
$ws.Range("M41").Value = -10835.75

$ws.Range("H45").Value = 2410.5
$ws.Range("I45").Value = 1032
$ws.Range("J45").Value = 3099.75
$ws.Range("K45").Value = 1032
$ws.Range("L45").Value = 3099.75
$ws.Range("M45").Value = -655
$ws.Range("N45").Value = -3853.75

$ws.Range("H110").Value = 2808.4827
$ws.Range("I110").Value = 3234.0527
$ws.Range("J110").Value = 1999.9
$ws.Range("K110").Value = 3234.0527
$ws.Range("L110").Value = 1999.9
$ws.Range("M110").Value = -1189.0527
$ws.Range("N110").Value = -6089.9

$ws.Range("H113").Value = 50000
$ws.Range("J113").Value = 50000
$ws.Range("L113").Value = 50000
$ws.Range("N113").Value = -58678

$ws.Range("H122").Value = 6669284.5
$ws.Range("I122").Value = 2313.081
$ws.Range("K122").Value = 6939.243
$ws.Range("M122").Value = -4489.243

$ws.Range("H132").Value = 17884690
$ws.Range("I132").Value = 2797.122
$ws.Range("K132").Value = 8391.366
$ws.Range("M132").Value = -5861.366

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 384
$ws.Range("I11").Value = 310.6842
$ws.Range("J11").Value = 848.3333
$ws.Range("K11").Value = 310.6842
$ws.Range("L11").Value = 848.3333
$ws.Range("M11").Value = -170.6842
$ws.Range("N11").Value = -1128.3333

$ws.Range("H61").Value = 80001
$ws.Range("J61").Value = 80001
$ws.Range("L61").Value = 80001
$ws.Range("N61").Value = -80627

$ws.Range("H107").Value = 2614.75
$ws.Range("I107").Value = 1286.4667
$ws.Range("K107").Value = 1286.4667
$ws.Range("M107").Value = 633.5333000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 807.75
$ws.Range("I16").Value = 829.06665
$ws.Range("J16").Value = 743.8
$ws.Range("K16").Value = 829.06665
$ws.Range("L16").Value = 743.8
$ws.Range("M16").Value = -542.06665
$ws.Range("N16").Value = -1317.8

$ws.Range("H87").Value = 75800.2
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52372

$ws.Range("H90").Value = 75800.2
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -161856

$ws.Range("H113").Value = 807.75
$ws.Range("I113").Value = 829.06665
$ws.Range("J113").Value = 743.8
$ws.Range("K113").Value = 829.06665
$ws.Range("L113").Value = 743.8
$ws.Range("M113").Value = 1340.93335
$ws.Range("N113").Value = -5083.8

$ws.Range("H134").Value = 2578.3333
$ws.Range("I134").Value = 2323.4736
$ws.Range("K134").Value = 6970.4208
$ws.Range("M134").Value = -4435.4208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 100000
$ws.Range("I14").Value = 100000
$ws.Range("K14").Value = 300000
$ws.Range("M14").Value = -299827

$ws.Range("H17").Value = 356.66666
$ws.Range("I17").Value = 240
$ws.Range("J17").Value = 380
$ws.Range("K17").Value = 720
$ws.Range("L17").Value = 1140
$ws.Range("M17").Value = -551
$ws.Range("N17").Value = -1478

$ws.Range("H125").Value = 111113130
$ws.Range("J125").Value = 166669170
$ws.Range("L125").Value = 500007510
$ws.Range("N125").Value = -500017350

$ws.Range("H126").Value = 119051290
$ws.Range("I126").Value = 125001416
$ws.Range("K126").Value = 375004248
$ws.Range("M126").Value = -374999308

$ws.Range("H131").Value = 6494734
$ws.Range("J131").Value = 1973.4445
$ws.Range("L131").Value = 5920.333500000001
$ws.Range("N131").Value = -16000.3335

$ws.Range("H132").Value = 5559925.5
$ws.Range("I132").Value = 1367.875
$ws.Range("K132").Value = 12310.875
$ws.Range("M132").Value = -9780.875

$ws.Range("H133").Value = 3033
$ws.Range("J133").Value = 3033
$ws.Range("L133").Value = 9099
$ws.Range("N133").Value = -19219

$ws.Range("H134").Value = 811.2857
$ws.Range("I134").Value = 811.2857
$ws.Range("K134").Value = 2433.8571
$ws.Range("M134").Value = 2636.1429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 779199.75
$ws.Range("I33").Value = 38400
$ws.Range("J33").Value = 1519999.5
$ws.Range("K33").Value = 38400
$ws.Range("L33").Value = 1519999.5
$ws.Range("M33").Value = -38148
$ws.Range("N33").Value = -1520503.5

$ws.Range("H119").Value = 40200
$ws.Range("J119").Value = 40200
$ws.Range("L119").Value = 40200
$ws.Range("N119").Value = -49876

$ws.Range("H132").Value = 4727.206
$ws.Range("I132").Value = 4840.1514
$ws.Range("K132").Value = 14520.4542
$ws.Range("M132").Value = -11990.4542

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2813.7273
$ws.Range("I7").Value = 2756.1304
$ws.Range("K7").Value = 2756.1304
$ws.Range("M7").Value = -2644.1304

$ws.Range("H16").Value = 2034.1818
$ws.Range("I16").Value = 2075.111
$ws.Range("J16").Value = 1850
$ws.Range("K16").Value = 2075.111
$ws.Range("L16").Value = 1850
$ws.Range("M16").Value = -1905.111
$ws.Range("N16").Value = -2190

$ws.Range("H40").Value = 1893.8667
$ws.Range("I40").Value = 1315.2142
$ws.Range("J40").Value = 9995
$ws.Range("K40").Value = 1315.2142
$ws.Range("L40").Value = 9995
$ws.Range("M40").Value = -1179.2142
$ws.Range("N40").Value = -10267

$ws.Range("H46").Value = 2088.8333

$ws.Range("H48").Value = 37495
$ws.Range("J48").Value = 37495
$ws.Range("L48").Value = 37495
$ws.Range("N48").Value = -38817

$ws.Range("H55").Value = 561.94116
$ws.Range("I55").Value = 601.75
$ws.Range("K55").Value = 601.75
$ws.Range("M55").Value = -428.75

$ws.Range("H100").Value = 4953.25
$ws.Range("I100").Value = 3592.7
$ws.Range("K100").Value = 3592.7
$ws.Range("M100").Value = -3051.7

$ws.Range("H108").Value = 60000
$ws.Range("J108").Value = 60000
$ws.Range("L108").Value = 60000
$ws.Range("N108").Value = -67680

$ws.Range("H122").Value = 4468518.5
$ws.Range("I122").Value = 3853.8948
$ws.Range("J122").Value = 13893922
$ws.Range("K122").Value = 11561.6844
$ws.Range("L122").Value = 41681766
$ws.Range("M122").Value = -9111.6844
$ws.Range("N122").Value = -41686666

$ws.Range("H126").Value = 2813.7273
$ws.Range("I126").Value = 2756.1304
$ws.Range("K126").Value = 8268.3912
$ws.Range("M126").Value = -5798.3912

$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960

$ws.Range("H136").Value = 1824527
$ws.Range("I136").Value = 2226977.5
$ws.Range("K136").Value = 6680932.5
$ws.Range("M136").Value = -6678382.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2174.2354
$ws.Range("I107").Value = 1107.5
$ws.Range("K107").Value = 3322.5
$ws.Range("M107").Value = -1402.5

$ws.Range("H136").Value = 2141.8262
$ws.Range("I136").Value = 2053.5
$ws.Range("J136").Value = 2459.8
$ws.Range("K136").Value = 6160.5
$ws.Range("L136").Value = 7379.400000000001
$ws.Range("M136").Value = -3610.5
$ws.Range("N136").Value = -12479.4

$ws.Range("H141").Value = 122998.336
$ws.Range("J141").Value = 129997.5
$ws.Range("L141").Value = 129997.5
$ws.Range("N141").Value = -140357.5
